$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 19 for the Bolivian Liga match;
# this shifts the former row 19 (Brazilian Serie A, Gremio vs Fluminense) down to row 20 automatically.
$ws.Rows.Item(19).Insert()

# Row 3: 4 changed cell(s)
$ws.Range('O3').Value2 = 1.35
$ws.Range('Q3').Value2 = 2
$ws.Range('S3').Value2 = 3.7
$ws.Range('AL3').Value2 = 80

# Row 4: 14 changed cell(s)
$ws.Range('F4').Value2 = 2.14
$ws.Range('H4').Value2 = 3.65
$ws.Range('I4').Value2 = 3.8
$ws.Range('L4').Value2 = 1.47
$ws.Range('N4').Value2 = 3.8
$ws.Range('O4').Value2 = 1.28
$ws.Range('P4').Value2 = 2.08
$ws.Range('Q4').Value2 = 1.9
$ws.Range('S4').Value2 = 1.9
$ws.Range('X4').Value2 = 17
$ws.Range('AB4').Value2 = 11
$ws.Range('AH4').Value2 = 970
$ws.Range('AJ4').Value2 = 970
$ws.Range('AN4').Value2 = 17.5

# Row 5: 3 changed cell(s)
$ws.Range('G5').Value2 = 2.32
$ws.Range('N5').Value2 = 2.98
$ws.Range('O5').Value2 = 1.4

# Row 6: 1 changed cell(s)
$ws.Range('H6').Value2 = 1.09

# Row 7: 1 changed cell(s)
$ws.Range('AL7').Value2 = 980

# Row 8: 1 changed cell(s)
$ws.Range('W8').Value2 = 1.72

# Row 9: 4 changed cell(s)
$ws.Range('M9').Value2 = 1.05
$ws.Range('P9').Value2 = 1.97
$ws.Range('Q9').Value2 = 1.7
$ws.Range('S9').Value2 = 2.78

# Row 10: 4 changed cell(s)
$ws.Range('L10').Value2 = 1.25
$ws.Range('N10').Value2 = 3.95
$ws.Range('V10').Value2 = 1.45
$ws.Range('W10').Value2 = 1.69

# Row 11: 38 changed cell(s)
$ws.Range('D11').Value2 = 'Bournemouth'
$ws.Range('E11').Value2 = 'Everton'
$ws.Range('F11').Value2 = 2.22
$ws.Range('G11').Value2 = 2.24
$ws.Range('H11').Value2 = 3.75
$ws.Range('I11').Value2 = 3.8
$ws.Range('J11').Value2 = 3.5
$ws.Range('K11').Value2 = 3.55
$ws.Range('L11').Value2 = 1.43
$ws.Range('M11').Value2 = 1.08
$ws.Range('N11').Value2 = 3.85
$ws.Range('O11').Value2 = 1.33
$ws.Range('P11').Value2 = 1.96
$ws.Range('Q11').Value2 = 2
$ws.Range('R11').Value2 = 1.37
$ws.Range('S11').Value2 = 3.55
$ws.Range('T11').Value2 = 1.8
$ws.Range('U11').Value2 = 2.2
$ws.Range('V11').Value2 = 1.35
$ws.Range('W11').Value2 = 1.81
$ws.Range('X11').Value2 = 13
$ws.Range('Y11').Value2 = 14.5
$ws.Range('Z11').Value2 = 26
$ws.Range('AA11').Value2 = 75
$ws.Range('AB11').Value2 = 9.800000000000001
$ws.Range('AC11').Value2 = 7.6
$ws.Range('AD11').Value2 = 15.5
$ws.Range('AE11').Value2 = 42
$ws.Range('AF11').Value2 = 13.5
$ws.Range('AG11').Value2 = 10.5
$ws.Range('AH11').Value2 = 17
$ws.Range('AI11').Value2 = 55
$ws.Range('AJ11').Value2 = 27
$ws.Range('AK11').Value2 = 23
$ws.Range('AL11').Value2 = 38
$ws.Range('AM11').Value2 = 95
$ws.Range('AN11').Value2 = 17
$ws.Range('AO11').Value2 = 44

# Row 12: 38 changed cell(s)
$ws.Range('D12').Value2 = 'Fulham'
$ws.Range('E12').Value2 = 'Man City'
$ws.Range('F12').Value2 = 5.3
$ws.Range('G12').Value2 = 5.4
$ws.Range('H12').Value2 = 1.71
$ws.Range('I12').Value2 = 1.72
$ws.Range('J12').Value2 = 4.4
$ws.Range('K12').Value2 = 4.5
$ws.Range('L12').Value2 = 1.32
$ws.Range('M12').Value2 = 1.05
$ws.Range('N12').Value2 = 5
$ws.Range('O12').Value2 = 1.23
$ws.Range('P12').Value2 = 2.34
$ws.Range('Q12').Value2 = 1.71
$ws.Range('R12').Value2 = 1.54
$ws.Range('S12').Value2 = 2.8
$ws.Range('T12').Value2 = 1.74
$ws.Range('U12').Value2 = 2.24
$ws.Range('V12').Value2 = 2.4
$ws.Range('W12').Value2 = 1.23
$ws.Range('X12').Value2 = 21
$ws.Range('Y12').Value2 = 10.5
$ws.Range('Z12').Value2 = 11.5
$ws.Range('AA12').Value2 = 17
$ws.Range('AB12').Value2 = 22
$ws.Range('AC12').Value2 = 9.6
$ws.Range('AD12').Value2 = 9.800000000000001
$ws.Range('AE12').Value2 = 15.5
$ws.Range('AF12').Value2 = 42
$ws.Range('AG12').Value2 = 20
$ws.Range('AH12').Value2 = 18.5
$ws.Range('AI12').Value2 = 28
$ws.Range('AJ12').Value2 = 130
$ws.Range('AK12').Value2 = 65
$ws.Range('AL12').Value2 = 60
$ws.Range('AM12').Value2 = 85
$ws.Range('AN12').Value2 = 55
$ws.Range('AO12').Value2 = 8

# Row 13: 11 changed cell(s)
$ws.Range('I13').Value2 = 4.2
$ws.Range('K13').Value2 = 4.5
$ws.Range('L13').Value2 = 1.19
$ws.Range('N13').Value2 = 7.2
$ws.Range('O13').Value2 = 1.14
$ws.Range('Q13').Value2 = 1.45
$ws.Range('R13').Value2 = 1.86
$ws.Range('S13').Value2 = 2.1
$ws.Range('V13').Value2 = 1.31
$ws.Range('W13').Value2 = 2.12
$ws.Range('AJ13').Value2 = 23

# Row 14: 8 changed cell(s)
$ws.Range('F14').Value2 = 1.75
$ws.Range('G14').Value2 = 1.76
$ws.Range('N14').Value2 = 4.7
$ws.Range('O14').Value2 = 1.25
$ws.Range('W14').Value2 = 2.3
$ws.Range('X14').Value2 = 19
$ws.Range('AC14').Value2 = 9.2
$ws.Range('AG14').Value2 = 9.8

# Row 15: 3 changed cell(s)
$ws.Range('K15').Value2 = 950
$ws.Range('N15').Value2 = 1.25
$ws.Range('S15').Value2 = 1.42

# Row 16: 18 changed cell(s)
$ws.Range('F16').Value2 = 1.97
$ws.Range('G16').Value2 = 2.06
$ws.Range('H16').Value2 = 4.4
$ws.Range('I16').Value2 = 5.1
$ws.Range('J16').Value2 = 3.25
$ws.Range('K16').Value2 = 3.55
$ws.Range('L16').Value2 = 1.49
$ws.Range('N16').Value2 = 3.25
$ws.Range('O16').Value2 = 1.4
$ws.Range('P16').Value2 = 1.76
$ws.Range('Q16').Value2 = 2.18
$ws.Range('S16').Value2 = 4.1
$ws.Range('U16').Value2 = 1.86
$ws.Range('V16').Value2 = 1.25
$ws.Range('W16').Value2 = 1.95
$ws.Range('AA16').Value2 = 140
$ws.Range('AE16').Value2 = 85
$ws.Range('AO16').Value2 = 120

# Row 18: 2 changed cell(s)
$ws.Range('I18').Value2 = 3.7
$ws.Range('Q18').Value2 = 1.87

# Row 19 (new): Bolivian Liga de Futbol Profesional - fill all cells
$ws.Range('A19').Value2 = 'Bolivian Liga de Futbol Profesional'
$ws.Range('B19').NumberFormat = '@'
$ws.Range('B19').Value2 = '2025-12-02'
$ws.Range('C19').Value2 = '20:30:00'
$ws.Range('D19').Value2 = 'Club Independiente Petrolero'
$ws.Range('E19').Value2 = 'Guabira'
$ws.Range('F19').Value2 = 1.04
$ws.Range('G19').Value2 = 980
$ws.Range('H19').Value2 = 1.04
$ws.Range('I19').Value2 = 980
$ws.Range('J19').Value2 = 1.02
$ws.Range('K19').Value2 = 980
$ws.Range('L19').Value2 = 1.01
$ws.Range('M19').Value2 = 1.01
$ws.Range('N19').Value2 = 1.25
$ws.Range('O19').Value2 = 1.25
$ws.Range('P19').Value2 = 1.25
$ws.Range('Q19').Value2 = 1.25
$ws.Range('R19').Value2 = 1.18
$ws.Range('S19').Value2 = 1.25
$ws.Range('T19').Value2 = 1.01
$ws.Range('U19').Value2 = 1.01
$ws.Range('V19').Value2 = 1.01
$ws.Range('W19').Value2 = 1.01
$ws.Range('X19').Value2 = 1000
$ws.Range('Y19').Value2 = 1000
$ws.Range('Z19').Value2 = 1000
$ws.Range('AA19').Value2 = 1000
$ws.Range('AB19').Value2 = 1000
$ws.Range('AC19').Value2 = 1000
$ws.Range('AD19').Value2 = 1000
$ws.Range('AE19').Value2 = 1000
$ws.Range('AF19').Value2 = 1000
$ws.Range('AG19').Value2 = 1000
$ws.Range('AH19').Value2 = 1000
$ws.Range('AI19').Value2 = 1000
$ws.Range('AJ19').Value2 = 1000
$ws.Range('AK19').Value2 = 1000
$ws.Range('AL19').Value2 = 1000
$ws.Range('AM19').Value2 = 1000
$ws.Range('AN19').Value2 = 1000
$ws.Range('AO19').Value2 = 1000
